$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Duel Decks Nissa vs. Ob Nixilis (DDR)',
    'Abyssal Gatekeeper',
    'Abyssal Specter',
    'Akroma, Angel of Wrath',
    'Angelic Benediction',
    'Angelic Page',
    'Angelic Protector',
    'Angel of Mercy',
    'Angel''s Feather',
    'Angelsong',
    'Barren Moor',
    'Barter in Blood',
    'Breeding Pit',
    'Cackling Imp',
    'Charging Paladin',
    'Consume Spirit',
    'Corrupt',
    'Cruel Edict',
    'Daggerclaw Imp',
    'Dark Banishing',
    'Dark Ritual',
    'Demonic Tutor',
    'Demon''s Horn',
    'Demon''s Jester',
    'Duress',
    'Dusk Imp',
    'Faith''s Fetters',
    'Fallen Angel',
    'Foul Imp',
    'Healing Salve',
    'Icatian Priest',
    'Kuro, Pitlord',
    'Lord of the Pit',
    'Luminous Angel',
    'Marble Diamond',
    'Oni Possession',
    'Otherworldly Journey',
    'Overeager Apprentice',
    'Pacifism',
    'Plains',
    'Plains',
    'Plains',
    'Plains',
    'Promise of Power',
    'Reiver Demon',
    'Reya Dawnbringer',
    'Righteous Cause',
    'Secluded Steppe',
    'Serra Advocate',
    'Serra Angel',
    'Serra''s Boon',
    'Serra''s Embrace',
    'Soot Imp',
    'Souldrinker',
    'Stinkweed Imp',
    'Sustainer of the Realm',
    'Swamp',
    'Swamp',
    'Swamp',
    'Swamp',
    'Twilight Shepherd',
    'Unholy Strength',
    'Venerable Monk'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
